$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update G2 value (269 -> 267), which ripples through dependent formulas
$ws.Range("G2").Value = 267

# 2. New formula cell Q25 = 267*20000, formatted like the "Rp" currency cells (style copied from I7)
$ws.Range("Q25").Formula = "=267*20000"
$ws.Range("I7").Copy() | Out-Null
$ws.Range("Q25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 3. New text cell P26 = "y"
$ws.Range("P26").Value = "y"

# 4. Set column Q width like the diff shows (width 12, bestFit/autofit)
$ws.Range("Q25").EntireColumn.AutoFit() | Out-Null

# 5. Update the view: scroll to top-left F1 and select Q25
$ws.Range("Q25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6
